# Update the "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Each entry is: row number -> new value for column F
$updates = @{
    4  = 1398
    5  = 470
    6  = 205
    10 = 299
    11 = 342
    12 = 339
    13 = 1819
    17 = 713
    20 = 4357
    22 = 314
    23 = 1174
    24 = 508
    25 = 49
    26 = 725
    28 = 377
    30 = 190
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
